# Reorder the "Periodo Mora" detail rows (B16:G24) grouping them by worker
# instead of by period, per commit: "Elimna EC anteriores y se agregan
# nuevos, se modifica base de datos".
#
# The set of underlying records is unchanged; only the row order (and,
# as a consequence, the shared-string table order) changes so that all
# periods for a given worker are listed together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row order for B16:G24 (Tipo Doc, N Doc, Nombre, Periodo, Valor Mora, Salario Basico)
$rows = @(
    @("CC", "73213301",  "JHON JAIRO BARRIOS SEPULVEDA",  "1912", 40000, 1000000),
    @("CC", "73213301",  "JHON JAIRO BARRIOS SEPULVEDA",  "1911", 40000, 1000000),
    @("CC", "73213301",  "JHON JAIRO BARRIOS SEPULVEDA",  "1910", 40000, 1000000),
    @("CC", "73213301",  "JHON JAIRO BARRIOS SEPULVEDA",  "1909", 40000, 1000000),
    @("CC", "1047403776","BELKIS MARIA FLOREZ GONZALEZ",  "1912", 33125, 878000),
    @("CC", "1047403776","BELKIS MARIA FLOREZ GONZALEZ",  "1911", 18771, 878000),
    @("CC", "1128048125","YULY PAULINA MUÑOZ OSPINO",     "1911", 33125, 743000),
    @("CC", "1128048125","YULY PAULINA MUÑOZ OSPINO",     "1910", 33125, 743000),
    @("CC", "1128048125","YULY PAULINA MUÑOZ OSPINO",     "1909", 33125, 743000)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 2).Value = $data[0]   # B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value = $data[1]   # C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value = $data[2]   # D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value = $data[3]   # E - Periodo Mora
    $ws.Cells.Item($r, 6).Value = $data[4]   # F - Valor Mora
    $ws.Cells.Item($r, 7).Value = $data[5]   # G - Salario Basico
}
